# Add season-record columns (Wins / Losses / Ties) to the team roster sheet.
# Every player row gets the same team season record: 93 wins, 69 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 should look like the rest of the header row
# (bold, centered, thin-bordered) - copy the formatting from the last
# existing header cell (AC1) before setting the new header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row.
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 93
    $ws.Cells.Item($row, 31).Value = 69
    $ws.Cells.Item($row, 32).Value = 0
}
